$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 7727.9375
$ws.Range("I70").Value = 2431.2632
$ws.Range("K70").Value = 7293.7896
$ws.Range("M70").Value = -7023.7896
$ws.Range("H73").Value = 7727.9375
$ws.Range("I73").Value = 2431.2632
$ws.Range("K73").Value = 7293.7896
$ws.Range("M73").Value = -6357.7896
$ws.Range("H112").Value = 1550596
$ws.Range("J112").Value = 1622602.8
$ws.Range("L112").Value = 4867808.4
$ws.Range("N112").Value = -4870024.4
$ws.Range("H138").Value = 5182.93
$ws.Range("J138").Value = 5908.8203
$ws.Range("L138").Value = 17726.4609
$ws.Range("N138").Value = -28006.4609

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9778.546
$ws.Range("I32").Value = 9403.276
$ws.Range("K32").Value = 9403.276
$ws.Range("M32").Value = -9116.276
$ws.Range("H61").Value = 2759.25
$ws.Range("I61").Value = 1856.2046
$ws.Range("K61").Value = 1856.2046
$ws.Range("M61").Value = -1644.2046
$ws.Range("H97").Value = 968.2105
$ws.Range("I97").Value = 611.58826
$ws.Range("J97").Value = 3999.5
$ws.Range("K97").Value = 611.58826
$ws.Range("L97").Value = 3999.5
$ws.Range("M97").Value = -115.58826
$ws.Range("N97").Value = -4991.5
$ws.Range("H110").Value = 15773.866
$ws.Range("I110").Value = 16079.9
$ws.Range("K110").Value = 16079.9
$ws.Range("M110").Value = -14034.9
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -49840
$ws.Range("H136").Value = 2759.25
$ws.Range("I136").Value = 1856.2046
$ws.Range("K136").Value = 5568.6138
$ws.Range("M136").Value = -3018.6138

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 9999
$ws.Range("I54").Value = 9999
$ws.Range("K54").Value = 9999
$ws.Range("M54").Value = -9515
$ws.Range("H86").Value = 5048.7144
$ws.Range("I86").Value = 4983.4287
$ws.Range("K86").Value = 4983.4287
$ws.Range("M86").Value = -3860.4287
$ws.Range("H89").Value = 5048.7144
$ws.Range("I89").Value = 4983.4287
$ws.Range("K89").Value = 24917.1435
$ws.Range("M89").Value = -19301.1435
$ws.Range("H134").Value = 2284.7437
$ws.Range("I134").Value = 1999.9688
$ws.Range("J134").Value = 2482.848
$ws.Range("K134").Value = 5999.9064
$ws.Range("L134").Value = 7448.544
$ws.Range("M134").Value = -3464.9064
$ws.Range("N134").Value = -12518.544

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 671.4545000000001
$ws.Range("I22").Value = 476.55554
$ws.Range("K22").Value = 476.55554
$ws.Range("M22").Value = -126.55554
$ws.Range("H62").Value = 22197.2
$ws.Range("I62").Value = 34997.332
$ws.Range("K62").Value = 34997.332
$ws.Range("M62").Value = -34373.332
$ws.Range("H65").Value = 22197.2
$ws.Range("I65").Value = 34997.332
$ws.Range("K65").Value = 174986.66
$ws.Range("M65").Value = -171866.66
$ws.Range("H86").Value = 38466340
$ws.Range("J86").Value = 9983.75
$ws.Range("L86").Value = 9983.75
$ws.Range("N86").Value = -12229.75
$ws.Range("H89").Value = 38466340
$ws.Range("J89").Value = 9983.75
$ws.Range("L89").Value = 49918.75
$ws.Range("N89").Value = -61150.75
$ws.Range("H94").Value = 1564.4
$ws.Range("I94").Value = 1645.625
$ws.Range("K94").Value = 1645.625
$ws.Range("M94").Value = -1194.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 770.36
$ws.Range("I107").Value = 768.2857
$ws.Range("J107").Value = 773
$ws.Range("K107").Value = 2304.8571
$ws.Range("L107").Value = 2319
$ws.Range("M107").Value = -384.8571000000002
$ws.Range("N107").Value = -6159
$ws.Range("H129").Value = 30703514
$ws.Range("I129").Value = 53031090
$ws.Range("J129").Value = 3096.25
$ws.Range("K129").Value = 159093270
$ws.Range("L129").Value = 9288.75
$ws.Range("M129").Value = -159088270
$ws.Range("N129").Value = -19288.75
$ws.Range("H136").Value = 1607.381
$ws.Range("I136").Value = 1527.75
$ws.Range("K136").Value = 4583.25
$ws.Range("M136").Value = 516.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 217147.92
$ws.Range("I2").Value = 371497.56
$ws.Range("J2").Value = 1058.4
$ws.Range("K2").Value = 371497.56
$ws.Range("L2").Value = 1058.4
$ws.Range("M2").Value = -371384.56
$ws.Range("N2").Value = -1284.4
$ws.Range("H102").Value = 2179.6667
$ws.Range("I102").Value = 1791.1428
$ws.Range("J102").Value = 3539.5
$ws.Range("K102").Value = 1791.1428
$ws.Range("L102").Value = 3539.5
$ws.Range("M102").Value = -169.1428000000001
$ws.Range("N102").Value = -6783.5
$ws.Range("H123").Value = 33395.2
$ws.Range("J123").Value = 33395.2
$ws.Range("L123").Value = 33395.2
$ws.Range("N123").Value = -38295.2
$ws.Range("H126").Value = 9114.414000000001
$ws.Range("I126").Value = 7286.143
$ws.Range("J126").Value = 10820.8
$ws.Range("K126").Value = 21858.429
$ws.Range("L126").Value = 32462.4
$ws.Range("M126").Value = -19388.429
$ws.Range("N126").Value = -37402.39999999999
$ws.Range("H132").Value = 54585.684
$ws.Range("I132").Value = 73095.07000000001
$ws.Range("J132").Value = 2759.4
$ws.Range("K132").Value = 219285.21
$ws.Range("L132").Value = 8278.200000000001
$ws.Range("M132").Value = -216755.21
$ws.Range("N132").Value = -13338.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 8026.8
$ws.Range("I55").Value = 385.41666
$ws.Range("K55").Value = 385.41666
$ws.Range("M55").Value = -212.41666
$ws.Range("H127").Value = 80238.336
$ws.Range("J127").Value = 80238.336
$ws.Range("L127").Value = 80238.336
$ws.Range("N127").Value = -90158.336
$ws.Range("H136").Value = 2365.5833
$ws.Range("I136").Value = 2064.4482
$ws.Range("J136").Value = 3613.1428
$ws.Range("K136").Value = 6193.344599999999
$ws.Range("L136").Value = 10839.4284
$ws.Range("M136").Value = -3643.344599999999
$ws.Range("N136").Value = -15939.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 216349.8
$ws.Range("I2").Value = 431850
$ws.Range("J2").Value = 162474.75
$ws.Range("K2").Value = 431850
$ws.Range("L2").Value = 162474.75
$ws.Range("M2").Value = -431738
$ws.Range("N2").Value = -162698.75
$ws.Range("H41").Value = 8865.286
$ws.Range("I41").Value = 8339.5
$ws.Range("K41").Value = 8339.5
$ws.Range("M41").Value = -7949.5
$ws.Range("H62").Value = 22738754
$ws.Range("I62").Value = 12497
$ws.Range("J62").Value = 27789032
$ws.Range("K62").Value = 12497
$ws.Range("L62").Value = 27789032
$ws.Range("M62").Value = -11873
$ws.Range("N62").Value = -27790280
$ws.Range("H65").Value = 22738754
$ws.Range("I65").Value = 12497
$ws.Range("J65").Value = 27789032
$ws.Range("K65").Value = 62485
$ws.Range("L65").Value = 138945160
$ws.Range("M65").Value = -59365
$ws.Range("N65").Value = -138951400
$ws.Range("H70").Value = 55998.5
$ws.Range("J70").Value = 49665
$ws.Range("L70").Value = 49665
$ws.Range("N70").Value = -50295
$ws.Range("H73").Value = 55998.5
$ws.Range("J73").Value = 49665
$ws.Range("L73").Value = 49665
$ws.Range("N73").Value = -51849
$ws.Range("H100").Value = 893.62164
$ws.Range("I100").Value = 703
$ws.Range("K100").Value = 1406
$ws.Range("M100").Value = -865
$ws.Range("H107").Value = 622.8823
$ws.Range("I107").Value = 704.7143
$ws.Range("K107").Value = 2114.1429
$ws.Range("M107").Value = -194.1428999999998
$ws.Range("H113").Value = 1951.7931
$ws.Range("I113").Value = 1955.7222
$ws.Range("K113").Value = 5867.1666
$ws.Range("M113").Value = -3697.1666
$ws.Range("H126").Value = 5438.92
$ws.Range("I126").Value = 4574.65
$ws.Range("J126").Value = 8896
$ws.Range("K126").Value = 13723.95
$ws.Range("L126").Value = 26688
$ws.Range("M126").Value = -11253.95
$ws.Range("N126").Value = -31628
